$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("students")

$numFmt = $ws.Range("C37").NumberFormat

# Row 38 (أحمد شوقي عبدالواسع محمد)
$ws.Range("C38").Value = 10
$ws.Range("D38").Value = 5
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 10
$ws.Range("G38").Value = 0
$ws.Range("C38:G38").NumberFormat = $numFmt
$ws.Range("J38").Value = "لا يوجد"

# Row 39 (أحمد محمد سعد محبوب)
$ws.Range("C39").Value = 10
$ws.Range("D39").Value = 5
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 10
$ws.Range("G39").Value = 1
$ws.Range("C39:G39").NumberFormat = $numFmt
$ws.Range("J39").Value = "لا يوجد"

# Update active selection to match the saved view state
$ws.Range("C40").Select()
